$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserData")

# Duplicate row 2 into a new row 3 (copies formatting + values)
$ws.Rows.Item(2).Copy() | Out-Null
$ws.Rows.Item(3).Insert() | Out-Null

# Update the new row's content for the second test user
$ws.Range("A3").Value = "JOHN1"
$ws.Range("B3").Value = "SMITH1"
$ws.Range("C3").Value = 352

# Add hyperlink for the email cell, matching row 2's hyperlink
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:kbisht@lendingpoint.com", "", "", "kbisht@lendingpoint.com") | Out-Null

# Move the active selection the way Excel would after these edits
$ws.Range("D4").Select() | Out-Null
